# Apply the account-balance export update to the "Export" worksheet.
#
# The sheet is a flat report of Conta / Nome / Saldo sorted descending by Saldo,
# followed by a blank row and a footer row with filter notes.
#
# This update:
#   - updates the balance for a few existing accounts
#   - adds a few brand-new account rows
#   - re-sorts the report rows (row 2 through the last data row) descending by Saldo,
#     which is how the rows end up re-ordered after the underlying values changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update balances on existing accounts (in place, before re-sort) ----
# Row 18: account 005142611 (Guilherme) 22242.97 -> 40243.2
$ws.Cells.Item(18, 3).Value2 = 40243.2

# Row 25: account 008153800 (Ricardo) 10121.93 -> 10121.95
$ws.Cells.Item(25, 3).Value2 = 10121.95

# Row 32: account 004584982 (Bruno) 2800.63 -> 6800.63
$ws.Cells.Item(32, 3).Value2 = 6800.63

# Row 51: account 002636063 (Leda) 578.17 -> 552.5
$ws.Cells.Item(51, 3).Value2 = 552.5

# ---- Insert brand-new account rows ----
# The last data row is row 462, followed by a blank row 463 and a footer row 464.
# Insert 3 fresh rows right before the blank row (i.e. before row 463), so the new
# records join the data block and push the blank/footer rows further down.
$ws.Rows.Item(463).Resize(3).Insert()

# Account numbers look numeric but must stay text (to keep their leading zeros),
# so force the "Conta" cells to text format before assigning their value.
$ws.Cells.Item(463, 1).NumberFormat = "@"
$ws.Cells.Item(463, 1).Value2 = "005142624"
$ws.Cells.Item(463, 2).Value2 = "Rodrigo"
$ws.Cells.Item(463, 3).Value2 = 13450

$ws.Cells.Item(464, 1).NumberFormat = "@"
$ws.Cells.Item(464, 1).Value2 = "008353082"
$ws.Cells.Item(464, 2).Value2 = "Pedro"
$ws.Cells.Item(464, 3).Value2 = 8147.35

$ws.Cells.Item(465, 1).NumberFormat = "@"
$ws.Cells.Item(465, 1).Value2 = "005186167"
$ws.Cells.Item(465, 2).Value2 = "Andrea"
$ws.Cells.Item(465, 3).Value2 = 2200

# The text-format trick above leaves the cells tagged with a "quote prefix" style;
# copy the plain formatting from an ordinary account-number cell on top of them so
# the new rows look just like the rest of the data (no leftover style/format).
$ws.Range("A2").Copy()
$ws.Range("A463:A465").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Re-sort the full data block (rows 2..465) descending by Saldo (column C) ----
$dataRange = $ws.Range("A2:C465")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("C2:C465"), 0, 2, 0, 0)
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = 0
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

Write-Host "Done applying Saldo export update"
